{"js": "// Optional rebate section for LED template.\n// - Wrap the rebate paragraph/sentence with <REBATE> ... </REBATE> tags so it\n//   can be switched on/off from database.json5.\n// - Rename the rebate-rate placeholder ${RR} -> ${ERR} (2 occurrences).\n// - Rename the payback placeholder ${PB} -> ${MPB}.\n// - Drop \"electric\" from \"through your electric utility company\".\n// - \"The total implementation cost is ${MIC}.\" becomes\n//   \"The modified implementation cost is ${MIC}.</REBATE>\" (closes the tag\n//   opened before the rebate explanation paragraph).\n\nconst body = context.document.body;\n\n// 1) Open the optional <REBATE> block and drop \"electric\" before \"utility company\".\nconst rebateIntro = body.search(\n  \"However, there could be energy efficiency rebates available through your electric utility company, which could potentially reduce the overall capital cost and thereby the payback period. The savings from the rebate is calculated below.\",\n  { matchCase: true }\n);\nrebateIntro.load(\"items\");\nawait context.sync();\n\nif (rebateIntro.items.length > 0) {\n  rebateIntro.items[0].insertText(\n    \"<REBATE>However, there could be energy efficiency rebates available through your utility company, which could potentially reduce the overall capital cost and thereby the payback period. The savings from the rebate is calculated below.\",\n    \"Replace\"\n  );\n}\n\n// 2) ${RR} -> ${ERR} (appears twice, in the RB formula rows).\nconst rrHits = body.search(\"${RR}\", { matchCase: true });\nrrHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < rrHits.items.length; i++) {\n  rrHits.items[i].insertText(\"${ERR}\", \"Replace\");\n}\n\n// 3) Close the optional <REBATE> block and mention the \"modified\" cost.\nconst micSentence = body.search(\"The total implementation cost is ${MIC}.\", {\n  matchCase: true,\n});\nmicSentence.load(\"items\");\nawait context.sync();\n\nif (micSentence.items.length > 0) {\n  micSentence.items[0].insertText(\n    \"The modified implementation cost is ${MIC}.</REBATE>\",\n    \"Replace\"\n  );\n}\n\n// 4) ${PB} -> ${MPB}.\nconst pbHits = body.search(\"${PB}\", { matchCase: true });\npbHits.load(\"items\");\nawait context.sync();\n\nif (pbHits.items.length > 0) {\n  pbHits.items[0].insertText(\"${MPB}\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Optional rebate for LED - can be switched from database.json5.\n#\n# - Wrap the rebate explanation with <REBATE> ... </REBATE> tags.\n# - Drop \"electric\" from \"through your electric utility company\".\n# - Rename the rebate-rate placeholder ${RR} -> ${ERR} (2 occurrences).\n# - \"The total implementation cost is ${MIC}.\" becomes\n#   \"The modified implementation cost is ${MIC}.</REBATE>\" (closing the tag).\n# - Rename the payback placeholder ${PB} -> ${MPB}.\n\n$d = $word.ActiveDocument\n\n# 1) Open the optional <REBATE> block and drop \"electric\" before \"utility company\".\n$findText1 = 'However, there could be energy efficiency rebates available through your electric utility company, which could potentially reduce the overall capital cost and thereby the payback period. The savings from the rebate is calculated below.'\n$replaceText1 = '<REBATE>However, there could be energy efficiency rebates available through your utility company, which could potentially reduce the overall capital cost and thereby the payback period. The savings from the rebate is calculated below.'\n$d.Content.Find.Execute($findText1, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText1, 2) | Out-Null\n\n# 2) ${RR} -> ${ERR} (appears twice, in the RB formula rows). Replace:=2 (wdReplaceAll) covers both.\n$d.Content.Find.Execute('${RR}', $false, $false, $false, $false, $false, $true, 1, $false, '${ERR}', 2) | Out-Null\n\n# 3) Close the optional <REBATE> block and mention the \"modified\" cost.\n$findText3 = 'The total implementation cost is ${MIC}.'\n$replaceText3 = 'The modified implementation cost is ${MIC}.</REBATE>'\n$d.Content.Find.Execute($findText3, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText3, 2) | Out-Null\n\n# 4) ${PB} -> ${MPB}.\n$d.Content.Find.Execute('${PB}', $false, $false, $false, $false, $false, $true, 1, $false, '${MPB}', 2) | Out-Null\n"}
